$d = $word.ActiveDocument

# Locate the end of the first bold run ("...XPS, R") and extend its text
# in-place so the split word "RTF" becomes contiguous again.
$rng = $d.Content
$found = $rng.Find.Execute("XPS, R", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
$rng.Collapse(0)  # wdCollapseEnd
$rng.InsertAfter("TF and TXT")

# Remove the now-duplicated old "TF and TXT" text that used to live in its
# own run right after the (soon to be removed) bookmark.
$rng2 = $d.Range($rng.End, $d.Content.End)
$found2 = $rng2.Find.Execute("TF and TXT", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
$rng2.Delete()

# Remove the obsolete "_GoBack" bookmark that used to split the run.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
